$d = $word.ActiveDocument

# Helper: find the paragraph whose text begins with a given prefix.
function Get-ParagraphByPrefix($doc, [string]$prefix) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $para = $doc.Paragraphs.Item($i)
        if ($para.Range.Text.StartsWith($prefix)) {
            return $para
        }
    }
    return $null
}

# ------------------------------------------------------------------
# 1) Gradient Boosting Machines (GBM) bullet:
#    extend the sentence with " - was used with Hidden Markov Model
#    (HMM)." and strike through the whole (now extended) sentence.
# ------------------------------------------------------------------
$gbmPara = Get-ParagraphByPrefix $d "Gradient Boosting Machines (GBM)"

$dash = [string][char]0x2013
$oldSentence = "Gradient Boosting Machines (GBM): An ensemble learning technique that builds models sequentially, each correcting the mistakes of the previous model."
$newSentence = "Gradient Boosting Machines (GBM): An ensemble learning technique that builds models sequentially, each correcting the mistakes of the previous model" + " " + $dash + " was used with Hidden Markov Model (HMM)."

$gbmPara.Range.Find.Execute($oldSentence, $true, $false, $false, $false, $false, $true, 1, $false, $newSentence, 2)
$gbmPara.Range.Font.StrikeThrough = 1

# ------------------------------------------------------------------
# 2) K-Nearest Neighbors (KNN) bullet: strike through entire sentence.
# ------------------------------------------------------------------
$knnPara = Get-ParagraphByPrefix $d "K-Nearest Neighbors (KNN)"
$knnPara.Range.Font.StrikeThrough = 1

# ------------------------------------------------------------------
# 3) Naive Bayes bullet: strike through entire sentence.
# ------------------------------------------------------------------
$nbPara = Get-ParagraphByPrefix $d "Naive Bayes"
$nbPara.Range.Font.StrikeThrough = 1

Write-Host "GBM updated : $($gbmPara.Range.Text)"
Write-Host "KNN updated : $($knnPara.Range.Text)"
Write-Host "NB  updated : $($nbPara.Range.Text)"
